# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" / "Valor Mora" table (rows 16-22) is re-sorted so the
# periods run in ascending order (2410, 2411, 2412, 2501, 2502, 2503, 2504)
# instead of the previous descending order (2504, 2503, 2502, 2501, 2412,
# 2411, 2410). The "Valor Mora" figures travel with their period, so the
# distinctive 39866 value (previously attached to period 2504 in row 16)
# now ends up attached to period 2504 in row 22, while every other period
# keeps its 52000 value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 16
$lastRow = 22

# Capture current "Periodo Mora" (E) and "Valor Mora" (F) pairs for the
# data rows before touching anything.
$periods = @()
$values = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $periods += , $ws.Cells.Item($r, 5).Value()
    $values += , $ws.Cells.Item($r, 6).Value()
}

# Re-write the rows in reverse order (ascending period), pairing each
# period back up with the value it originally had.
$rowCount = $lastRow - $firstRow + 1
for ($i = 0; $i -lt $rowCount; $i++) {
    $targetRow = $firstRow + $i
    $sourceIndex = $rowCount - 1 - $i
    $ws.Cells.Item($targetRow, 5).Value = $periods[$sourceIndex]
    $ws.Cells.Item($targetRow, 6).Value = $values[$sourceIndex]
}

$wb.Save()
